# Add a "time_taken" metadata column (F) to the panel worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cell F1 --------------------------------------------------
# Give it the same look as the other header cells (bold, boxed border,
# centered) by copying the formatting straight off the neighboring
# header cell E1 ("panel").
$ws.Range("F1").Value = "time_taken"
$ws.Range("F1").Font.Bold = $ws.Range("E1").Font.Bold
$ws.Range("F1").HorizontalAlignment = $ws.Range("E1").HorizontalAlignment
$ws.Range("F1").VerticalAlignment = $ws.Range("E1").VerticalAlignment
$ws.Range("F1").Borders.LineStyle = $ws.Range("E1").Borders.LineStyle

# --- Data rows F2:F20 --------------------------------------------------
# Per-row "time_taken" timestamps recorded for each gene/panel entry.
$timestamps = @(
    "2021-10-05 10:50:14.785266",
    "2021-10-05 10:50:14.785277",
    "2021-10-05 10:50:14.785281",
    "2021-10-05 10:50:14.785284",
    "2021-10-05 10:50:14.785288",
    "2021-10-05 10:50:14.785291",
    "2021-10-05 10:50:14.785294",
    "2021-10-05 10:50:14.785297",
    "2021-10-05 10:50:14.785301",
    "2021-10-05 10:50:14.785304",
    "2021-10-05 10:50:14.785307",
    "2021-10-05 10:50:14.785310",
    "2021-10-05 10:50:14.785313",
    "2021-10-05 10:50:14.785316",
    "2021-10-05 10:50:14.785319",
    "2021-10-05 10:50:14.785322",
    "2021-10-05 10:50:14.785326",
    "2021-10-05 10:50:14.785329",
    "2021-10-05 10:50:14.785332"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}
